# Updated cryptos list on Tue Aug  6 02:48:05 UTC 2024 with GitHub Actions
# Refreshes per-row Price (col D) and Volume(1h) (col E) figures, and
# swaps the VeChain/Bittensor rows (48/49) to reflect their new rank order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $text) {
    # Excel auto-detects plain decimal strings (e.g. "0.997") as numbers,
    # which would change the cell's stored type away from text. Prefixing
    # with a literal apostrophe forces text entry (same as typing it by
    # hand); the ClearFormats() afterwards drops the resulting quote-prefix
    # marker so the cell format matches the untouched cells again.
    $cell = $ws.Range($addr)
    if ($text -match '^[+-]?[0-9]*\.?[0-9]+$') {
        $cell.Value = "'" + $text
        $cell.ClearFormats()
    } else {
        $cell.Value = $text
    }
}

Set-TextCell 'D2' '55.464.93'
Set-TextCell 'E2' '  +3.00%  '
Set-TextCell 'D3' '2.480.50'
Set-TextCell 'E3' '  +10.11%  '
Set-TextCell 'D5' '480.86'
Set-TextCell 'E5' '  +8.47%  '
Set-TextCell 'D6' '139.11'
Set-TextCell 'E6' '  +9.16%  '
Set-TextCell 'D7' '0.997'
Set-TextCell 'E7' '  -0.12%  '
Set-TextCell 'D8' '0.509'
Set-TextCell 'E8' '  +8.32%  '
Set-TextCell 'D9' '2.506.52'
Set-TextCell 'E9' '  +11.29%  '
Set-TextCell 'D10' '0.0983'
Set-TextCell 'E10' '  +8.07%  '
Set-TextCell 'D11' '5.44'
Set-TextCell 'E11' '  +1.82%  '
Set-TextCell 'D12' '0.325'
Set-TextCell 'E12' '  +5.78%  '
Set-TextCell 'E13' '  +0.00%  '
Set-TextCell 'D14' '2.940.50'
Set-TextCell 'E14' '  +11.12%  '
Set-TextCell 'D15' '55.579.56'
Set-TextCell 'E15' '  +3.09%  '
Set-TextCell 'D16' '20.38'
Set-TextCell 'E16' '  +9.04%  '
Set-TextCell 'D17' '0.0000136'
Set-TextCell 'E17' '  +14.72%  '
Set-TextCell 'D18' '2.518.76'
Set-TextCell 'E18' '  +11.07%  '
Set-TextCell 'D19' '4.31'
Set-TextCell 'E19' '  +6.30%  '
Set-TextCell 'D20' '320.39'
Set-TextCell 'E20' '  +6.62%  '
Set-TextCell 'D21' '9.84'
Set-TextCell 'E21' '  +4.73%  '
Set-TextCell 'D22' '0.997'
Set-TextCell 'E22' '  -0.13%  '
Set-TextCell 'D23' '5.65'
Set-TextCell 'E23' '  +7.12%  '
Set-TextCell 'D24' '57.68'
Set-TextCell 'E24' '  +3.87%  '
Set-TextCell 'E25' '  +3.95%  '
Set-TextCell 'D26' '0.164'
Set-TextCell 'E26' '  +3.80%  '
Set-TextCell 'D27' '0.401'
Set-TextCell 'E27' '  +8.25%  '
Set-TextCell 'D28' '2.622.09'
Set-TextCell 'E28' '  +11.26%  '
Set-TextCell 'D29' '7.37'
Set-TextCell 'E29' '  +8.75%  '
Set-TextCell 'D30' '0.0₃0765'
Set-TextCell 'E30' '  +9.38%  '
Set-TextCell 'D31' '0.998'
Set-TextCell 'E31' '  +0.19%  '
Set-TextCell 'D32' '149.74'
Set-TextCell 'E32' '  +4.37%  '
Set-TextCell 'D33' '18.06'
Set-TextCell 'E33' '  +7.70%  '
Set-TextCell 'E34' '  +10.51%  '
Set-TextCell 'D35' '5.17'
Set-TextCell 'E35' '  +9.58%  '
Set-TextCell 'D36' '3.68'
Set-TextCell 'E36' '  +2.01%  '
Set-TextCell 'E37' '  +10.46%  '
Set-TextCell 'D38' '0.839'
Set-TextCell 'E38' '  +0.83%  '
Set-TextCell 'E39' '  +4.99%  '
Set-TextCell 'D40' '0.613'
Set-TextCell 'E40' '  +21.35%  '
Set-TextCell 'D41' '0.993'
Set-TextCell 'E41' '  -0.24%  '
Set-TextCell 'D42' '1.32'
Set-TextCell 'E42' '  +7.16%  '
Set-TextCell 'D43' '0.0543'
Set-TextCell 'E43' '  +10.05%  '
Set-TextCell 'D44' '3.36'
Set-TextCell 'E44' '  +7.07%  '
Set-TextCell 'D45' '1.982.27'
Set-TextCell 'E45' '  +3.64%  '
Set-TextCell 'D46' '10.13'
Set-TextCell 'E46' '  -1.72%  '
Set-TextCell 'D47' '0.0896'
Set-TextCell 'E47' '  +11.30%  '
Set-TextCell 'B48' 'Bittensor'
Set-TextCell 'C48' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextCell 'D48' '249.97'
Set-TextCell 'E48' '  +33.50%  '
Set-TextCell 'B49' 'VeChain'
Set-TextCell 'C49' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell 'D49' '0.0221'
Set-TextCell 'E49' '  +8.06%  '
Set-TextCell 'D50' '17.54'
Set-TextCell 'E50' '  +8.38%  '
Set-TextCell 'D51' '4.38'
Set-TextCell 'E51' '  +9.47%  '
